$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Header: grade -> level
$ws.Range("E1").Value = "level"

# Row 2
$ws.Range("E2").Value = "المستوى الأول"

# Row 3
$ws.Range("E3").Value = "المستوى الأول"

# Row 4
$ws.Range("B4").Value = "08:00"
$ws.Range("C4").Value = "10:00"
$ws.Range("D4").Value = "الرياضيات"
$ws.Range("E4").Value = "المستوى الثاني"

# Row 5
$ws.Range("E5").Value = "المستوى الأول"

# Row 6
$ws.Range("E6").Value = "المستوى الثاني"

# Row 7
$ws.Range("B7").Value = "08:00"
$ws.Range("C7").Value = "10:00"
$ws.Range("D7").Value = "العلوم"
$ws.Range("E7").Value = "المستوى الثاني"

# Row 8
$ws.Range("E8").Value = "المستوى الأول"
# F8 must stay a text cell ("2"), not become a number - force text format,
# then revert to the original (unstyled) look so no stray style sticks around.
$ws.Range("F8").NumberFormat = "@"
$ws.Range("F8").Value = "2"
$ws.Range("F8").Style = "Normal"

# Row 9
$ws.Range("E9").Value = "المستوى الثاني"

# Row 10
$ws.Range("E10").Value = "المستوى الأول"

# Row 11
$ws.Range("E11").Value = "المستوى الثاني"
